# Apply the diff: add missing H/J/M/O/Q values to rows 162-164 and
# append four new data rows (165-168) to the "Kennzahlen" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in previously-missing cells for existing rows 162-164 ---
$ws.Range("H162").Value = -1.285714285714286
$ws.Range("J162").Value = -1.538802220125937
$ws.Range("M162").Value = 6
$ws.Range("O162").Value = 0.1428571428571428
$ws.Range("Q162").Value = 7.142857142857143

$ws.Range("H163").Value = -1.285714285714286
$ws.Range("J163").Value = -1.538802220125937
$ws.Range("M163").Value = 6
$ws.Range("O163").Value = 0.1428571428571428
$ws.Range("Q163").Value = 7.142857142857143

$ws.Range("H164").Value = 0.1428571428571428
$ws.Range("J164").Value = 0.4673706193802354
$ws.Range("M164").Value = 7.857142857142857
$ws.Range("O164").Value = 0.1428571428571428
$ws.Range("Q164").Value = 7.571428571428571

# --- Append new rows 165-168 ---

# Row 165
$ws.Range("A164").Copy()
$ws.Range("A165:A168").PasteSpecial(-4122)

$ws.Range("A165").Value = 44054
$ws.Range("B165").Value = 2181
$ws.Range("C165").Value = 1070
$ws.Range("D165").Value = 101
$ws.Range("E165").Value = 2008
$ws.Range("F165").Value = 72
$ws.Range("G165").Value = 1
$ws.Range("H165").Value = 0.1428571428571428
$ws.Range("I165").Value = 1.408450704225352
$ws.Range("J165").Value = 0.4673706193802354
$ws.Range("K165").Value = 11
$ws.Range("L165").Value = 2
$ws.Range("M165").Value = 7.857142857142857
$ws.Range("N165").Value = 1
$ws.Range("O165").Value = 0.1428571428571428
$ws.Range("P165").Value = 9
$ws.Range("Q165").Value = 7.571428571428571
$ws.Range("R165").Value = 392.6440009721584
$ws.Range("S165").Value = 432.5329452663918
$ws.Range("T165").Value = 6.872018756568841
$ws.Range("U165").Value = 7.561232480894385

# Row 166
$ws.Range("A166").Value = 44055
$ws.Range("B166").Value = 2181
$ws.Range("C166").Value = 1070
$ws.Range("D166").Value = 101
$ws.Range("E166").Value = 2008
$ws.Range("F166").Value = 72
$ws.Range("G166").Value = 0
$ws.Range("I166").Value = 0
$ws.Range("K166").Value = 0
$ws.Range("L166").Value = 0
$ws.Range("N166").Value = 0
$ws.Range("P166").Value = 0
$ws.Range("R166").Value = 392.6440009721584
$ws.Range("S166").Value = 432.5329452663918
$ws.Range("T166").Value = 6.872018756568841
$ws.Range("U166").Value = 7.561232480894385

# Row 167
$ws.Range("A167").Value = 44056
$ws.Range("B167").Value = 2212
$ws.Range("C167").Value = 1088
$ws.Range("D167").Value = 101
$ws.Range("E167").Value = 2028
$ws.Range("F167").Value = 83
$ws.Range("G167").Value = 11
$ws.Range("I167").Value = 15.27777777777778
$ws.Range("K167").Value = 31
$ws.Range("L167").Value = 18
$ws.Range("N167").Value = 0
$ws.Range("P167").Value = 20
$ws.Range("R167").Value = 398.2249106604377
$ws.Range("S167").Value = 439.8092004204058
$ws.Range("T167").Value = 11.31861912846633
$ws.Range("U167").Value = 9.901613963075981

# Row 168
$ws.Range("A168").Value = 44057
$ws.Range("B168").Value = 2212
$ws.Range("C168").Value = 1088
$ws.Range("D168").Value = 101
$ws.Range("E168").Value = 2028
$ws.Range("F168").Value = 83
$ws.Range("G168").Value = 0
$ws.Range("I168").Value = 0
$ws.Range("K168").Value = 0
$ws.Range("L168").Value = 0
$ws.Range("N168").Value = 0
$ws.Range("P168").Value = 0
$ws.Range("R168").Value = 398.2249106604377
$ws.Range("S168").Value = 439.8092004204058
$ws.Range("T168").Value = 11.31861912846633
$ws.Range("U168").Value = 9.901613963075981
